$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forms")

# --- Header row: new columns AF1:AK1 ---------------------------------------
# Written column-by-column (header then data cell below) so that new shared
# strings are interned in the same order as the target workbook.

# A13 is written first (new row, column A) before any header additions.
$ws.Range("A13").Value = "ProductRegistration"

# AF column (SKUitemNumber / numeric SKU)
$ws.Range("AF1").Value = "SKUitemNumber"
$ws.Range("AF1").Interior.Color = $ws.Range("A1").Interior.Color

# AG column (SerialNumber / numeric serial)
$ws.Range("AG1").Value = "SerialNumber"
$ws.Range("AG1").Interior.Color = $ws.Range("A1").Interior.Color

# AH column (ManufactureDate / text "0808")
$ws.Range("AH1").Value = "ManufactureDate"
$ws.Range("AH1").Interior.Color = $ws.Range("A1").Interior.Color

$ws.Range("AH13").Value = "'0808"
$ws.Range("AH13").Style = $ws.Range("M9").Style

# AI column (PurchasedAt / text "OXO Website")
$ws.Range("AI1").Value = "PurchasedAt"
$ws.Range("AI1").Interior.Color = $ws.Range("A1").Interior.Color

$ws.Range("AI13").Value = "OXO Website"

# AJ column (Price / numeric)
$ws.Range("AJ1").Value = "Price"
$ws.Range("AJ1").Interior.Color = $ws.Range("A1").Interior.Color

# AK1: same fill/style as the rest of the header row, no text
$ws.Range("AK1").Interior.Color = $ws.Range("A1").Interior.Color

# --- Re-style the old "consolas" header block (R1:AE1) ----------------------
# These used to carry a one-off Consolas font + yellow fill; restore them to
# the plain default-font + yellow-fill look used by the rest of row 1.
$ws.Range("R1:AE1").ClearFormats()
$ws.Range("R1:AE1").Interior.Color = $ws.Range("A1").Interior.Color

# --- New data row 13 (rest of the cells) ------------------------------------
$ws.Range("F13").Value = "test"
$ws.Range("G13").Value = "qa"

$ws.Range("H13").Value = "qatesting.lotuswave@gmail.com"
$ws.Hyperlinks.Add($ws.Range("H13"), "mailto:qatesting.lotuswave@gmail.com") | Out-Null

$ws.Range("J13").Value = "844 N colony rd"
$ws.Range("K13").Value = "Wallingford"
$ws.Range("L13").Value = "Connecticut"

$ws.Range("M13").Value = "'06492"
$ws.Range("M13").Style = $ws.Range("M9").Style

$ws.Range("N13").Value = "'9898989898"
$ws.Range("N13").Style = $ws.Range("N9").Style

$ws.Range("O13").Value = "product"
$ws.Range("P13").Value = "United States"

$ws.Range("Q13").Value = "'9/16/22"
$ws.Range("Q13").Style = $ws.Range("Q9").Style

$ws.Range("AF13").Value = 8718800
$ws.Range("AG13").Value = 23456789
$ws.Range("AJ13").Value = 199.99

# --- Sheet activation: "Forms" becomes the selected/visible tab ------------
$ws.Range("N15").Select()
$ws.Activate()
